$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was MITS5501/10/Lewis/F2F) -> becomes old Row 6 data (MITS4003/4/Tom/F2F)
$ws.Range("B2").Value = "8:00 AM to 10:00 AM"
$ws.Range("C2").Value = "MITS4003"
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = "Tom"
$ws.Range("F2").Value = "F2F"

# Row 3 (was MITS5503/12/Mike/F2F) -> becomes old Row 2 data (MITS5501/10/Lewis/F2F)
$ws.Range("B3").Value = "8:00 AM to 10:00 AM"
$ws.Range("C3").Value = "MITS5501"
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = "Lewis"
$ws.Range("F3").Value = "F2F"

# Row 4 (was MITS4001/1/Jim/Online) -> becomes old Row 3 data (MITS5503/12/Mike/F2F)
$ws.Range("B4").Value = "8:00 AM to 9:00 AM"
$ws.Range("C4").Value = "MITS5503"
$ws.Range("D4").Value = 12
$ws.Range("E4").Value = "Mike"
$ws.Range("F4").Value = "F2F"

# Row 6 (was MITS4003/4/Tom/F2F) -> becomes old Row 4 data (MITS4001/1/Jim/Online)
$ws.Range("B6").Value = "8:00 AM to 9:00 AM"
$ws.Range("C6").Value = "MITS4001"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "Jim"
$ws.Range("F6").Value = "Online"

# Row 10 and Row 11 are swapped
$ws.Range("C10").Value = "MITS5502"
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = "Jake"

$ws.Range("C11").Value = "MITS5004"
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = "Sean"

# Row 15 and Row 16 are swapped
$ws.Range("B15").Value = "2:00 PM to 3:00 PM"
$ws.Range("C15").Value = "MITS6500"
$ws.Range("D15").Value = 21
$ws.Range("E15").Value = "Keno"

$ws.Range("B16").Value = "2:00 PM to 4:00 PM"
$ws.Range("C16").Value = "MITS5003"
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = "Jay"
